$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("A12").Value = "011"
$ws.Range("B12").Value = 1344
$ws.Range("C12").Value = 729
$ws.Range("D12").Value = 1717
$ws.Range("E12").Value = 829
$ws.Range("F12").Value = "kick off button"

# Row 13
$ws.Range("A13").Value = "012"
$ws.Range("B13").Value = 746
$ws.Range("C13").Value = 924
$ws.Range("D13").Value = 1092
$ws.Range("E13").Value = 1012
$ws.Range("F13").Value = "go to scenario list"

# Copy styles from the row above (row 11) to keep formatting consistent
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A13:F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply values in case PasteSpecial(formats) touched anything, and ensure text format on column A
$ws.Range("A12").Value = "011"
$ws.Range("A13").Value = "012"

$ws.Range("F13").Select()
